# Refresh cryptocurrency price (col D) and 1h volume-change (col E) figures
# on Sheet1, matching the latest coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.163.02'
$ws.Range('E2').Value = '  +0.17%  '
$ws.Range('D3').Value = '1.826.12'
$ws.Range('E3').Value = '  -0.29%  '
$ws.Range('D4').Value = '''0.9985'
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('D5').Value = '''241.65'
$ws.Range('E5').Value = '  -0.67%  '
$ws.Range('D6').Value = '''0.6219'
$ws.Range('E6').Value = '  -0.72%  '
$ws.Range('D7').Value = '''0.9999'
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('D8').Value = '''0.07352'
$ws.Range('E8').Value = '  -1.91%  '
$ws.Range('E9').Value = '  -1.05%  '
$ws.Range('D10').Value = '''23.03'
$ws.Range('E10').Value = '  -1.24%  '
$ws.Range('D11').Value = '''0.07676'
$ws.Range('E11').Value = '  -0.37%  '
$ws.Range('D12').Value = '1.825.65'
$ws.Range('E12').Value = '  -0.11%  '
$ws.Range('D13').Value = '''4.961'
$ws.Range('E13').Value = '  -1.16%  '
$ws.Range('D14').Value = '''0.6653'
$ws.Range('E14').Value = '  -0.46%  '
$ws.Range('D15').Value = '''82.39'
$ws.Range('E15').Value = '  -0.46%  '
$ws.Range('D16').Value = '''0.000008998'
$ws.Range('E16').Value = '  -4.17%  '
$ws.Range('D17').Value = '''5.849'
$ws.Range('E17').Value = '  -2.27%  '
$ws.Range('D18').Value = '29.138.39'
$ws.Range('E18').Value = '  +0.18%  '
$ws.Range('D19').Value = '2.070.59'
$ws.Range('E19').Value = '  +1.32%  '
$ws.Range('D20').Value = '''238.53'
$ws.Range('E20').Value = '  +6.80%  '
$ws.Range('D21').Value = '''12.43'
$ws.Range('E21').Value = '  -1.34%  '
$ws.Range('D22').Value = '''0.9996'
$ws.Range('E22').Value = '  -0.21%  '
$ws.Range('D23').Value = '''7.258'
$ws.Range('E23').Value = '  +1.72%  '
$ws.Range('D24').Value = '''0.9997'
$ws.Range('E24').Value = '  -0.28%  '
$ws.Range('D25').Value = '''158.11'
$ws.Range('E25').Value = '  -1.20%  '
$ws.Range('E26').Value = '  +2.18%  '
$ws.Range('D27').Value = '''8.486'
$ws.Range('E27').Value = '  -0.22%  '
$ws.Range('D28').Value = '''17.67'
$ws.Range('E28').Value = '  -1.34%  '
$ws.Range('D29').Value = '''1.485'
$ws.Range('E29').Value = '  -0.45%  '
$ws.Range('D30').Value = '''0.05576'
$ws.Range('E30').Value = '  -3.57%  '
$ws.Range('D31').Value = '''4.090'
$ws.Range('E31').Value = '  -0.81%  '
$ws.Range('D32').Value = '''4.096'
$ws.Range('E32').Value = '  -1.55%  '
$ws.Range('E33').Value = '  -0.46%  '
$ws.Range('D34').Value = '''1.843'
$ws.Range('E34').Value = '  +0.71%  '
$ws.Range('D35').Value = '''0.7339'
$ws.Range('E35').Value = '  -1.11%  '
$ws.Range('E36').Value = '  -0.35%  '
$ws.Range('D37').Value = '''2.627'
$ws.Range('E37').Value = '  -1.67%  '
$ws.Range('D38').Value = '''2.840'
$ws.Range('E38').Value = '  +2.70%  '
$ws.Range('D39').Value = '1.213.07'
$ws.Range('E39').Value = '  -1.57%  '
$ws.Range('D40').Value = '''0.01766'
$ws.Range('E40').Value = '  -0.42%  '
$ws.Range('D41').Value = '''6.308'
$ws.Range('E41').Value = '  -2.98%  '
$ws.Range('D42').Value = '''0.9124'
$ws.Range('E42').Value = '  +2.20%  '
$ws.Range('D43').Value = '''1.000'
$ws.Range('E43').Value = '  -0.06%  '
$ws.Range('D44').Value = '''101.68'
$ws.Range('E44').Value = '  -0.55%  '
$ws.Range('D45').Value = '1.973.48'
$ws.Range('D46').Value = '''64.69'
$ws.Range('E46').Value = '  -1.80%  '
$ws.Range('D47').Value = '''0.5084'
$ws.Range('E47').Value = '  -0.24%  '
$ws.Range('E48').Value = '  -2.87%  '
$ws.Range('D49').Value = '''0.4025'
$ws.Range('E49').Value = '  -1.06%  '
$ws.Range('D50').Value = '''9.123'
$ws.Range('E50').Value = '  +1.13%  '
$ws.Range('D51').Value = '''0.05760'
$ws.Range('E51').Value = '  -1.17%  '
